# Apply the "MHD" ALE/GAGE/AGAGE instrument combination sheet edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet Sheet1 -> MHD
$ws.Name = "MHD"

# --- Column widths (A:I default width, B:C / D widened to fit header text) ---
$ws.Columns("A:A").ColumnWidth = 9.998697916666666
$ws.Columns("B:C").ColumnWidth = 14.998697916666666
$ws.Columns("D:D").ColumnWidth = 14.666666666666666
$ws.Columns("E:I").ColumnWidth = 9.998697916666666

# --- Header / comment rows (A1:A7) ---
$ws.Range("A1").Value = "# File specifying when to use the various ALE/GAGE/AGAGE instruments"
$ws.Range("A2").Value = "# A blank entry in either the start or end date means that the time is unbounded at that side"
$ws.Range("A3").Value = "# If there is no entry for a species it is assumed that it is only measured on the GCMS-Medusa"
$ws.Range("A4").Value = "# Time periods can overlap if you want to keep both instruments"
$ws.Range("A5").Value = '# An "x" indicates that this instrument should not be included in timeseries for this species'
$ws.Range("A6").Value = "# Date format should be YYYY-MM-DD HH:MM"
$ws.Range("A7").Value = '# NOTE: ENSURE CELLS ARE FORMATTED AS TEXT, NOT DATES. IF PASTING VALUES ENSURE YOU "MATCH DESTINATION FORMATTING" TO PREVENT EXCEL FROM CONVERTING TO DATE AND TIME'
$ws.Range("A1:A7").NumberFormat = "@"

# --- Table header row (row 8), bold ---
$ws.Range("A8").Value = "Species"
$ws.Range("B8").Value = "GCMD start"
$ws.Range("C8").Value = "GCMD end"
$ws.Range("D8").Value = "Picarro start"
$ws.Range("E8").Value = "Picarro end"
$ws.Range("A8:E8").NumberFormat = "@"
$ws.Range("A8:E8").Font.Bold = $true

# --- View tweaks ---
$ws.Range("A10").Select()
